# Integration test vouchers: add "Voucher Recommendation" sheet data,
# and tweak a few sheet views/selections that Excel re-saved.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Voucher Recommendation" sheet after "My Voucher" ---
$afterSheet = $wb.Worksheets.Item("My Voucher")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "Voucher Recommendation"

# Header row
$newSheet.Range("A1").Value = "testCase"
$newSheet.Range("B1").Value = "transactionId"
$newSheet.Range("C1").Value = "result"

# Row 2
$newSheet.Range("A2").Value = "Valid parameters"
$newSheet.Range("B2").Value = "'true"
$newSheet.Range("C2").Value = '{"id":4,"name":"Cashback Rp 2.500 buat beli pulsa","voucherTypeName":"cashback","discount":0,"maxDeduction":0,"value":2500,"filePath":"https://res.cloudinary.com/darwmcfjo/image/upload/v1591548482/WhatsApp_Image_2020-05-30_at_7.27.48_PM_nessez.jpg","expiryDate":1279584000000}'

# Row 3
$newSheet.Range("A3").Value = "Have no vouchers"
$newSheet.Range("B3").Value = "'true"
$newSheet.Range("C3").Value = "[]"

# Column widths (matches the sheet as saved by Excel)
$newSheet.Columns.Item(1).ColumnWidth = 15.33203125
$newSheet.Columns.Item(2).ColumnWidth = 11.83203125
$newSheet.Columns.Item(3).ColumnWidth = 254.33203125

# Selection / active cell on the new sheet
$newSheet.Range("B7").Select()

# --- 2. sheet1 ("Voucher Promotion"): selection changed to a block range ---
$wsPromotion = $wb.Worksheets.Item("Voucher Promotion")
$wsPromotion.Range("A1:C2").Select()

# --- 3. sheet2 ("My Voucher"): no longer the tab-selected sheet; new active cell ---
$wsMyVoucher = $wb.Worksheets.Item("My Voucher")
$wsMyVoucher.Range("C13").Select()

# --- 4. sheet4 ("Voucher Details"): column A width shrunk; new active cell ---
$wsDetails = $wb.Worksheets.Item("Voucher Details")
$wsDetails.Columns.Item(1).ColumnWidth = 9
$wsDetails.Range("C9").Select()

# Make the new sheet the active / tab-selected sheet, matching activeTab="2"
$newSheet.Activate()
